# Atualiza o arquivo modelo de importação de transações (modelo-transacoes.xlsx)
# Substitui os dados de exemplo (3 linhas) por um conjunto mais completo de
# exemplos de Receita/Despesa (12 linhas), nas colunas:
#   A = Data, B = Descrição, C = Valor, D = Tipo, E = Categoria

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# As datas são textos (ex.: "2024-09-01"), não números de data do Excel.
# O prefixo de apóstrofo força o Excel a gravar o valor como texto literal,
# preservando o mesmo tipo de célula (string) usado no arquivo original.

# Linha 2
$ws.Cells.Item(2, 1).Value = "'2024-09-01"
$ws.Cells.Item(2, 2).Value = "Venda de produtos no varejo"
$ws.Cells.Item(2, 3).Value = 1250
$ws.Cells.Item(2, 4).Value = "Receita"
$ws.Cells.Item(2, 5).Value = "Varejo"

# Linha 3
$ws.Cells.Item(3, 1).Value = "'2024-09-02"
$ws.Cells.Item(3, 2).Value = "Venda em atacado para distribuidor"
$ws.Cells.Item(3, 3).Value = 3500
$ws.Cells.Item(3, 4).Value = "Receita"
$ws.Cells.Item(3, 5).Value = "Atacado"

# Linha 4
$ws.Cells.Item(4, 1).Value = "'2024-09-03"
$ws.Cells.Item(4, 2).Value = "Rendimento de investimentos"
$ws.Cells.Item(4, 3).Value = 450
$ws.Cells.Item(4, 4).Value = "Receita"
$ws.Cells.Item(4, 5).Value = "Investimentos"

# Linha 5
$ws.Cells.Item(5, 1).Value = "'2024-09-04"
$ws.Cells.Item(5, 2).Value = "Serviços de consultoria"
$ws.Cells.Item(5, 3).Value = 800
$ws.Cells.Item(5, 4).Value = "Receita"
$ws.Cells.Item(5, 5).Value = "Outros"

# Linha 6
$ws.Cells.Item(6, 1).Value = "'2024-09-05"
$ws.Cells.Item(6, 2).Value = "Venda online de produtos"
$ws.Cells.Item(6, 3).Value = 680
$ws.Cells.Item(6, 4).Value = "Receita"
$ws.Cells.Item(6, 5).Value = "Varejo"

# Linha 7
$ws.Cells.Item(7, 1).Value = "'2024-09-01"
$ws.Cells.Item(7, 2).Value = "Aluguel do escritório"
$ws.Cells.Item(7, 3).Value = 1200
$ws.Cells.Item(7, 4).Value = "Despesa"
$ws.Cells.Item(7, 5).Value = "Fixo"

# Linha 8
$ws.Cells.Item(8, 1).Value = "'2024-09-02"
$ws.Cells.Item(8, 2).Value = "Compra de matéria-prima"
$ws.Cells.Item(8, 3).Value = 850
$ws.Cells.Item(8, 4).Value = "Despesa"
$ws.Cells.Item(8, 5).Value = "Variável"

# Linha 9
$ws.Cells.Item(9, 1).Value = "'2024-09-03"
$ws.Cells.Item(9, 2).Value = "Compra de equipamentos"
$ws.Cells.Item(9, 3).Value = 2500
$ws.Cells.Item(9, 4).Value = "Despesa"
$ws.Cells.Item(9, 5).Value = "Investimento"

# Linha 10
$ws.Cells.Item(10, 1).Value = "'2024-09-04"
$ws.Cells.Item(10, 2).Value = "Campanha de marketing digital"
$ws.Cells.Item(10, 3).Value = 300
$ws.Cells.Item(10, 4).Value = "Despesa"
$ws.Cells.Item(10, 5).Value = "Mkt"

# Linha 11
$ws.Cells.Item(11, 1).Value = "'2024-09-05"
$ws.Cells.Item(11, 2).Value = "Despesas administrativas"
$ws.Cells.Item(11, 3).Value = 150
$ws.Cells.Item(11, 4).Value = "Despesa"
$ws.Cells.Item(11, 5).Value = "Outros"

# Linha 12
$ws.Cells.Item(12, 1).Value = "'2024-09-06"
$ws.Cells.Item(12, 2).Value = "Salários dos funcionários"
$ws.Cells.Item(12, 3).Value = 3200
$ws.Cells.Item(12, 4).Value = "Despesa"
$ws.Cells.Item(12, 5).Value = "Fixo"

# Linha 13
$ws.Cells.Item(13, 1).Value = "'2024-09-07"
$ws.Cells.Item(13, 2).Value = "Combustível para entrega"
$ws.Cells.Item(13, 3).Value = 180
$ws.Cells.Item(13, 4).Value = "Despesa"
$ws.Cells.Item(13, 5).Value = "Variável"

# Garante que o intervalo de erros ignorados (número armazenado como texto)
# acompanhe a nova área de dados, caso a propriedade seja suportada.
try {
    $ws.Range("A1:E13").Errors.Item(6).Ignore = $true
} catch {
}

Write-Host "Modelo de transacoes atualizado (A1:E13)."
